$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A2" = 2
    "A3" = 1
    "A4" = 1
    "A5" = 1
    "A6" = 1
    "A11" = 2
    "A15" = 1
    "A17" = 1
    "A21" = 2
    "A22" = 2
    "A23" = 2
    "A26" = 1
    "A31" = 1
    "A33" = 1
    "A36" = 1
    "A40" = 1
    "A45" = 2
    "A47" = 2
    "A48" = 1
    "A50" = 1
    "A52" = 1
    "A54" = 2
    "A55" = 2
    "A56" = 1
    "A59" = 1
    "A61" = 2
    "A62" = 1
    "A66" = 2
    "A67" = 1
    "A70" = 2
    "A71" = 1
    "A73" = 2
    "A74" = 2
    "A76" = 1
    "A77" = 1
    "A80" = 1
    "A81" = 1
    "A82" = 1
    "A87" = 2
    "A89" = 2
    "A91" = 2
    "A94" = 1
    "A95" = 2
    "A96" = 1
    "A97" = 2
    "A98" = 1
    "A103" = 1
    "A105" = 2
    "A107" = 2
    "A108" = 1
    "A110" = 1
    "A111" = 1
    "A112" = 1
    "A113" = 1
    "A114" = 1
    "A118" = 1
    "A119" = 1
    "A120" = 1
    "A122" = 2
    "A123" = 1
    "A127" = 2
    "A128" = 2
    "A129" = 1
    "A130" = 2
    "A132" = 1
    "A134" = 1
    "A142" = 2
    "A145" = 1
    "A146" = 1
    "A151" = 1
    "A154" = 1
    "A162" = 2
    "A165" = 2
    "A167" = 1
    "A168" = 1
    "A169" = 1
    "A170" = 2
    "A175" = 1
    "A177" = 1
    "A179" = 1
    "A181" = 1
    "A184" = 2
    "A187" = 2
    "A189" = 1
    "A195" = 2
    "A197" = 1
    "A199" = 1
    "A200" = 1
    "A201" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

